$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Level"
$ws.Range("B1").Value = "Respiratory HI"
$ws.Range("C1").Value = "Liver HI"
$ws.Range("D1").Value = "Neurological HI"
$ws.Range("E1").Value = "Developmental HI"
$ws.Range("F1").Value = "Reproductive HI"
$ws.Range("G1").Value = "Kidney HI"
$ws.Range("H1").Value = "Ocular HI"
$ws.Range("I1").Value = "Endocrine HI"
$ws.Range("J1").Value = "Hematological HI"
$ws.Range("K1").Value = "Immunological HI"
$ws.Range("L1").Value = "Skeletal HI"
$ws.Range("M1").Value = "Spleen HI"
$ws.Range("N1").Value = "Thyroid HI"
$ws.Range("O1").Value = "Whole body HI"
